# Add a new "Templates" worksheet at the end of the workbook for the
# new Word Mail Merge Template functionality.
$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Templates"

# Header row (same yellow header fill used on the other sheets)
$ws.Range("A1").Value = "title"
$ws.Range("B1").Value = "description"
$ws.Range("C1").Value = "version"
$ws.Range("D1").Value = "tags"
$ws.Range("A1:D1").Interior.Color = 65535

# Sample template rows
$ws.Range("A2").Value = "Template1"
$ws.Range("B2").Value = "Template1 description"
$ws.Range("C2").Value = "v1"
$ws.Range("D2").Value = "Template 1 tags"

$ws.Range("A3").Value = "Template2"
$ws.Range("B3").Value = "Template2 description"
$ws.Range("C3").Value = "v2"
$ws.Range("D3").Value = "Template 2 tags"

# Column widths to fit the content (mirrors the other sheets' best-fit widths)
$ws.Columns.Item(1).ColumnWidth = 8.830729166666666
$ws.Columns.Item(2).ColumnWidth = 18.498697916666668
$ws.Columns.Item(3).ColumnWidth = 6.053385416666667
$ws.Columns.Item(4).ColumnWidth = 13.166666666666666

# Leave the selection / active tab on the new sheet, matching the edit
$ws.Range("D3").Select() | Out-Null
